$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$c = $s.Shapes.AddConnector(2, 100,100,100,300)
$members = $c | Get-Member
Write-Output $members
